$d = $word.ActiveDocument

# --- Fix "protec" -> "protect" ------------------------------------------
# Replacing a range that spans from inside the preceding run into the
# misspelled word collapses the (now corrected) text into a single run,
# which also drops the w:proofErr spell-check markers that wrapped it.
$full = $d.Content.Text
$anchor = "to protec "
$idx = $full.IndexOf($anchor)
$r = $d.Range($idx, $idx + $anchor.Length)
$r.Text = "to protect "

# --- Fix "balck" -> "black" ----------------------------------------------
$full2 = $d.Content.Text
$anchor2 = "the balck kids"
$idx2 = $full2.IndexOf($anchor2)
$r2 = $d.Range($idx2, $idx2 + $anchor2.Length)
$r2.Text = "the black kids"

# --- Re-split the sentence into the run layout from the target ----------
# Both fixes above merged everything from "protect" onward into one run.
# Re-create the expected run boundaries by toggling a character property
# on (and immediately back off) each sub-range; Word breaks a new run at
# the edges of the touched range even though the formatting ends up
# identical to its neighbours.
$full3 = $d.Content.Text
$sentenceStart = $full3.IndexOf("He sent the us army to protect and let the black kids enter school.")

function Split-Run([int]$startOffset, [int]$endOffset) {
    $pos1 = $sentenceStart + $startOffset
    $pos2 = $sentenceStart + $endOffset
    $rr = $d.Range($pos1, $pos2)
    $rr.Font.Bold = 1
    $rr.Font.Bold = 0
}

# Desired runs (offsets into the sentence above):
#   [0, 23)  "He sent the us army to "   <- left untouched
#   [23, 30) "protect"
#   [30, 44) " and let the b"
#   [44, 45) "l"
#   [45, 67) "ack kids enter school."
Split-Run 23 30
Split-Run 30 44
Split-Run 44 45
